$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 14 (pushes existing rows 14-36 down to 15-37)
$ws.Rows("14:14").Insert()

# Populate the new row 14 with a new weekly data point.
# Values mirror the entry that used to sit in row 14, except for the
# date (one week later) and the reported volume.
$ws.Range("A14").Value = 9
$ws.Range("B14").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C14").Value = "Metropolitana"
$ws.Range("D14").Value = 44953
$ws.Range("E14").Value = 13
$ws.Range("F14").Value = 100112010
$ws.Range("G14").Value = "Achicoria"
$ws.Range("H14").Value = "Sin especificar"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 90
$ws.Range("K14").Value = 7000
$ws.Range("L14").Value = 7000
$ws.Range("M14").Value = 7000
$ws.Range("N14").Value = "$/caja 16 unidades"
$ws.Range("O14").Value = "Provincia de Quillota"
$ws.Range("P14").Value = 438
$ws.Range("Q14").Value = 16
$ws.Range("R14").Value = "Hortaliza"
